$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13; this shifts the existing rows 13..75
# down to 14..76 and grows the used range to A1:R76.
$ws.Rows(13).Insert()

# Populate the newly inserted row 13 with the new record.
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Vega Monumental Concepción"
$ws.Range("C13").Value = "Bíobío"
$ws.Range("D13").Value = 44600
$ws.Range("D13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E13").Value = 8
$ws.Range("F13").Value = 100112024
$ws.Range("G13").Value = "Choclo"
$ws.Range("H13").Value = "Choclero"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 6000
$ws.Range("K13").Value = 150
$ws.Range("L13").Value = 200
$ws.Range("M13").Value = 175
$ws.Range("N13").Value = "$/unidad"
$ws.Range("O13").Value = "Región del Maule"
$ws.Range("P13").Value = 175
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = "Hortaliza"
